# The edit corresponds to moving the data-row that was in sheet row 15
# down to become the last of the four rows 15-18 (rows 16, 17 and 18
# each shift up by one row). Row 19 and below, and rows 14 and above,
# are untouched.
#
# We reproduce that row-move with whole-row Range.Copy() operations
# (which preserve each cell's original type - numbers stay numbers,
# text that merely looks like a date/time stays text - instead of
# Range.Value assignment, which would risk Excel re-interpreting a
# literal string like "2026-01-31" as a date serial number).
#
# Because this engine's Copy()-paste only overwrites cells that have
# content in the source range (it does not blank out a destination
# cell when the corresponding source cell is empty), each destination
# row is explicitly cleared immediately before it receives the copied
# row, so columns that should end up empty really end up empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fullRow = "A{0}:AY{0}"
$stageRow = 200   # scratch row well outside the used range (A1:AY20)

# 1) Stage the original row 15 out of the way before it gets overwritten.
$ws.Range(($fullRow -f $stageRow)).ClearContents() | Out-Null
$ws.Range(($fullRow -f 15)).Copy($ws.Range(($fullRow -f $stageRow)))

# 2) Shift rows 16 -> 15, 17 -> 16, 18 -> 17.
$ws.Range(($fullRow -f 15)).ClearContents() | Out-Null
$ws.Range(($fullRow -f 16)).Copy($ws.Range(($fullRow -f 15)))

$ws.Range(($fullRow -f 16)).ClearContents() | Out-Null
$ws.Range(($fullRow -f 17)).Copy($ws.Range(($fullRow -f 16)))

$ws.Range(($fullRow -f 17)).ClearContents() | Out-Null
$ws.Range(($fullRow -f 18)).Copy($ws.Range(($fullRow -f 17)))

# 3) Drop the staged original row 15 content into row 18.
$ws.Range(($fullRow -f 18)).ClearContents() | Out-Null
$ws.Range(($fullRow -f $stageRow)).Copy($ws.Range(($fullRow -f 18)))

# 4) Clean up the scratch row so it doesn't linger in the saved file.
$ws.Range(($fullRow -f $stageRow)).ClearContents() | Out-Null

Write-Output "Row 15 moved to position 18 (rows 16-18 shifted up one)."
